# Delete row 672 ("「とても愛しています」..." post) so that all subsequent
# rows shift up by one (row 673 becomes 672, ..., row 876 becomes 875).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(672).Delete()
